# "All email endpoints done"
# Update the endpoint documentation sheet so the email & user-login rows use
# the new Request/Response-suffixed DTO type names, then move the active
# selection and widen column H to fit the new (longer) body-type text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: user / POST /api/users (Login) ---------------------------------
# Body type:   UserLoginData      -> UserLoginDataRequest
# Return type: User                -> UserResponse
$ws.Cells.Item(3, 8).Value2 = "UserLoginDataRequest"
$ws.Cells.Item(3, 9).Value2 = "UserResponse"

# --- Rows 9-13: user / GET find-all* endpoints ------------------------------
# Return type: List<User> -> List<UserResponse>
$ws.Cells.Item(9, 9).Value2  = "List<UserResponse>"
$ws.Cells.Item(10, 9).Value2 = "List<UserResponse>"
$ws.Cells.Item(11, 9).Value2 = "List<UserResponse>"
$ws.Cells.Item(12, 9).Value2 = "List<UserResponse>"
$ws.Cells.Item(13, 9).Value2 = "List<UserResponse>"

# --- Row 17: email / POST /api/emails (add) --------------------------------
# Body type:   String -> EmailRequest
# Return type: Email  -> EmailResponse
$ws.Cells.Item(17, 8).Value2 = "EmailRequest"
$ws.Cells.Item(17, 9).Value2 = "EmailResponse"

# --- Row 18: email / DELETE /api/emails (delete) ---------------------------
# Body type: String -> EmailRequest
$ws.Cells.Item(18, 8).Value2 = "EmailRequest"

# --- Sheet view / formatting tweaks -----------------------------------------
# The author's selection ended on H4 instead of D14 when they saved.
$ws.Range("H4").Select()

# Column H ("Body") needs to be widened to fit "UserLoginDataRequest" /
# "EmailRequest" now that it's no longer just "String".
$ws.Columns.Item(8).ColumnWidth = 20.67
